$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "45.710.13"
Set-TextCell $ws.Range("E2") "  +7.07%  "
Set-TextCell $ws.Range("D3") "2.404.46"
Set-TextCell $ws.Range("E3") "  +4.52%  "
Set-TextCell $ws.Range("E4") "  -0.01%  "
Set-TextCell $ws.Range("D5") "115.11"
Set-TextCell $ws.Range("E5") "  +9.86%  "
Set-TextCell $ws.Range("D6") "317.81"
Set-TextCell $ws.Range("E6") "  +2.02%  "
Set-TextCell $ws.Range("D7") "0.639"
Set-TextCell $ws.Range("E7") "  +2.30%  "
Set-TextCell $ws.Range("E8") "  -0.18%  "
Set-TextCell $ws.Range("D9") "0.630"
Set-TextCell $ws.Range("E9") "  +4.13%  "
Set-TextCell $ws.Range("D10") "42.72"
Set-TextCell $ws.Range("E10") "  +7.38%  "
Set-TextCell $ws.Range("E11") "  +3.37%  "
Set-TextCell $ws.Range("D12") "8.78"
Set-TextCell $ws.Range("E12") "  +6.20%  "
Set-TextCell $ws.Range("E13") "  +2.66%  "
Set-TextCell $ws.Range("D14") "1.01"
Set-TextCell $ws.Range("E14") "  +2.89%  "
Set-TextCell $ws.Range("D15") "15.99"
Set-TextCell $ws.Range("E15") "  +4.15%  "
Set-TextCell $ws.Range("D16") "2.773.77"
Set-TextCell $ws.Range("E16") "  -0.23%  "
Set-TextCell $ws.Range("D17") "2.412.95"
Set-TextCell $ws.Range("E17") "  +5.01%  "
Set-TextCell $ws.Range("D18") "45.736.80"
Set-TextCell $ws.Range("E18") "  +6.79%  "
Set-TextCell $ws.Range("D19") "7.55"
Set-TextCell $ws.Range("E19") "  +3.38%  "
Set-TextCell $ws.Range("E20") "  +3.88%  "
Set-TextCell $ws.Range("D21") "13.56"
Set-TextCell $ws.Range("E21") "  -0.37%  "
Set-TextCell $ws.Range("D22") "75.03"
Set-TextCell $ws.Range("E22") "  +2.24%  "
Set-TextCell $ws.Range("D23") "3.58"
Set-TextCell $ws.Range("E23") "  +3.51%  "
Set-TextCell $ws.Range("D24") "264.96"
Set-TextCell $ws.Range("E24") "  -1.42%  "
Set-TextCell $ws.Range("E25") "  +7.50%  "
Set-TextCell $ws.Range("E26") "  -0.78%  "
Set-TextCell $ws.Range("D27") "11.44"
Set-TextCell $ws.Range("E27") "  +5.47%  "
Set-TextCell $ws.Range("D28") "7.62"
Set-TextCell $ws.Range("E28") "  +6.30%  "
Set-TextCell $ws.Range("D29") "40.89"
Set-TextCell $ws.Range("E29") "  +13.02%  "
Set-TextCell $ws.Range("E30") "  +2.29%  "
Set-TextCell $ws.Range("E31") "  +16.27%  "
Set-TextCell $ws.Range("D32") "22.80"
Set-TextCell $ws.Range("E32") "  +2.17%  "
Set-TextCell $ws.Range("D33") "173.20"
Set-TextCell $ws.Range("E33") "  +5.09%  "
Set-TextCell $ws.Range("D34") "2.95"
Set-TextCell $ws.Range("E34") "  +12.86%  "
Set-TextCell $ws.Range("B35") "Stellar"
Set-TextCell $ws.Range("C35") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D35") "0.133"
Set-TextCell $ws.Range("E35") "  +2.08%  "
Set-TextCell $ws.Range("B36") "RenderToken"
Set-TextCell $ws.Range("C36") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D36") "5.02"
Set-TextCell $ws.Range("E36") "  +10.03%  "
Set-TextCell $ws.Range("E37") "  +7.52%  "
Set-TextCell $ws.Range("D38") "4.30"
Set-TextCell $ws.Range("E38") "  +18.81%  "
Set-TextCell $ws.Range("D39") "3.13"
Set-TextCell $ws.Range("E39") "  +11.96%  "
Set-TextCell $ws.Range("E40") "  +5.72%  "
Set-TextCell $ws.Range("E41") "  +12.71%  "
Set-TextCell $ws.Range("E42") "  +12.20%  "
Set-TextCell $ws.Range("E43") "  +6.25%  "
Set-TextCell $ws.Range("B44") "BitcoinSV"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell $ws.Range("D44") "99.43"
Set-TextCell $ws.Range("E44") "  -10.17%  "
Set-TextCell $ws.Range("B45") "MultiversX"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws.Range("D45") "72.11"
Set-TextCell $ws.Range("E45") "  +1.60%  "
Set-TextCell $ws.Range("D46") "86.18"
Set-TextCell $ws.Range("E46") "  +11.54%  "
Set-TextCell $ws.Range("E47") "  -0.46%  "
Set-TextCell $ws.Range("D48") "5.84"
Set-TextCell $ws.Range("E48") "  +13.72%  "
Set-TextCell $ws.Range("D49") "116.71"
Set-TextCell $ws.Range("E49") "  +5.44%  "
Set-TextCell $ws.Range("D50") "9.62"
Set-TextCell $ws.Range("E50") "  +11.42%  "
Set-TextCell $ws.Range("D51") "1.61"
Set-TextCell $ws.Range("E51") "  +12.39%  "
